$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new "2022" column (S), mirroring the layout/format
# already used by the adjacent "2021" column (R).
$ws.Range("R3:R33").Copy($ws.Range("S3:S33"))

# Header year
$ws.Range("S3").Value = 2022

# Men (rows 4-6), Women (rows 7-9) ... repeated per-region blocks
$ws.Range("S4").Value = 5.5
$ws.Range("S5").Value = 8.5
$ws.Range("S6").Value = 2.6

$ws.Range("S7").Value = 16.3
$ws.Range("S8").Value = 25.2
$ws.Range("S9").Value = 7.1

$ws.Range("S10").Value = 1.6
$ws.Range("S11").Value = 3.2
$ws.Range("S12").Value = "-"

$ws.Range("S13").Value = 7.5
$ws.Range("S14").Value = 10.5
$ws.Range("S15").Value = 4.5

$ws.Range("S16").Value = 11.4
$ws.Range("S17").Value = 16.1
$ws.Range("S18").Value = 6.6

$ws.Range("S19").Value = 1.2
$ws.Range("S20").Value = 2.1
$ws.Range("S21").Value = 0.3

$ws.Range("S22").Value = 1.5
$ws.Range("S23").Value = 2.9
$ws.Range("S24").Value = 0

$ws.Range("S25").Value = 0.9
$ws.Range("S26").Value = 1.7
$ws.Range("S27").Value = 0.2

$ws.Range("S28").Value = 14.3
$ws.Range("S29").Value = 22.7
$ws.Range("S30").Value = 7.3

$ws.Range("S31").Value = 1.1
$ws.Range("S32").Value = 2.2
$ws.Range("S33").Value = "-"

# Move the selection like the source workbook ends up (next empty header cell)
$ws.Range("T3").Select()
